$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" under the duplicate_image_filename column (E) for every
# stimulus/data row (rows 2-21), per the commit message:
# "add the NA's under duplicate_image_filename"
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
